$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.038.48"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").Value = "2.365.12"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'500.70"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").Value = "'128.47"
$ws.Range("E6").Value = "  -3.51%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D9").Value = "2.368.79"
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("D10").Value = "'0.0978"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "'4.75"
$ws.Range("E12").Value = "  +3.38%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "2.783.24"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "55.994.95"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "'21.40"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "2.418.88"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "'9.98"
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("D20").Value = "'4.03"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "'305.82"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("D22").Value = "'6.27"
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "'65.49"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("E26").Value = "  -3.58%  "
$ws.Range("E27").Value = "  -5.64%  "
$ws.Range("E28").Value = "  -4.92%  "
$ws.Range("D29").Value = "'171.99"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").Value = "0.0₃0710"
$ws.Range("E30").Value = "  -2.96%  "
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").Value = "'5.73"
$ws.Range("E34").Value = "  -7.01%  "
$ws.Range("D35").Value = "'1.08"
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("D36").Value = "'17.61"
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("E37").Value = "  -5.32%  "
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("D39").Value = "'36.03"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Value = "'1.37"
$ws.Range("E41").Value = "  -6.09%  "
$ws.Range("D42").Value = "'129.58"
$ws.Range("E42").Value = "  -4.03%  "
$ws.Range("D43").Value = "'3.35"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E44").Value = "  -6.41%  "
$ws.Range("D45").Value = "'0.561"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").Value = "'0.0901"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").Value = "'239.53"
$ws.Range("E47").Value = "  -6.60%  "
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("E49").Value = "  -3.77%  "
$ws.Range("D50").Value = "'17.04"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  -0.66%  "
